# Update specific statistics cells in the single results table.
# Table layout: col 1 = row label, col 2 = Danish, col 3 = Dutch,
# col 4 = English, col 5 = French, col 6 = Spanish, col 7 = Total.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellValue($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.End = $r.End - 1   # drop the end-of-cell marker
    $r.Text = $newText
}

# (mean) SLAXIMP block
Set-CellValue $t 5  5 "365"     # Mean - French: 368 -> 365
Set-CellValue $t 6  5 "361"     # Median - French: 373 -> 361
Set-CellValue $t 10 5 "42"      # Number of nonmissing values - French: 43 -> 42
Set-CellValue $t 10 7 "363"     # Number of nonmissing values - Total: 364 -> 363

# Enslaved person mortality rate block
Set-CellValue $t 17 5 "41"      # Number of nonmissing values - French: 42 -> 41
Set-CellValue $t 17 7 "362"     # Number of nonmissing values - Total: 363 -> 362

# Total net expenditure in g. of silver per enslaved person block
Set-CellValue $t 19 5 "3,152"   # Mean - French: 3,153 -> 3,152
Set-CellValue $t 19 7 "2,643"   # Mean - Total: 2,645 -> 2,643
Set-CellValue $t 20 5 "2,678"   # Median - French: 2,713 -> 2,678
Set-CellValue $t 21 5 "1,829"   # Standard deviation - French: 1,807 -> 1,829
Set-CellValue $t 21 7 "1,493"   # Standard deviation - Total: 1,491 -> 1,493
Set-CellValue $t 24 5 "42"      # Number of nonmissing values - French: 43 -> 42
Set-CellValue $t 24 7 "363"     # Number of nonmissing values - Total: 364 -> 363
